$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: a corrected "Coco" / "Royal Park" entry (same show as row 2, but
# with the right date/time; columns B-D intentionally left blank)
$ws.Range("A5").Value = "Coco"
$ws.Range("E5").Value = "Royal Park"
$ws.Range("F5").Value = "Dec 4"
$ws.Range("G5").Value = "1:10pm"

# New row 6: a new "Shark Tank" episode entry
$ws.Range("A6").Value = "Shark Tank"
$ws.Range("B6").Value = "Episode 307"
$ws.Range("C6").Value = "(Season 3, Episode 2)"
$ws.Range("D6").Value = "Ideas include a training system for salespeople, a management system for cargo trucks, family-friendly Las Vegas entertainment and a cat-portrait business. Also: a follow-up on a Season 2 potty-training product for cats."
$ws.Range("E6").Value = "CNBC"
$ws.Range("F6").Value = "Dec 5"
$ws.Range("G6").Value = "7:00pm"
